$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells I1 and J1, styled like existing header row (copy style from H1)
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I (I0) and J (IF), rows 2-19
$values = @{
    2  = @(5, 6)
    3  = @(7, 7)
    4  = @(7, 8)
    5  = @(7, 7)
    6  = @(9, 9)
    7  = @(10, 10)
    8  = @(10, 10)
    9  = @(8, 8)
    10 = @(9, 9)
    11 = @(9, 9)
    12 = @(9, 9)
    13 = @(9, 9)
    14 = @(3, 3)
    15 = @(9, 9)
    16 = @(9, 9)
    17 = @(9, 9)
    18 = @(9, 9)
    19 = @(5, 5)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
